$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "desc" row (it goes where "name" used to be) by
# shifting the "name"/"sity" row data down one row first.
$ws.Range("A7").Value = $ws.Range("A6").Value
$ws.Range("B7").Value = $ws.Range("B6").Value
$ws.Range("A6").Value = $ws.Range("A5").Value
$ws.Range("B6").Value = $ws.Range("B5").Value

# A7 is a brand-new cell with no formatting yet - give it the same
# bold/bordered/centered label style used by the rest of column A.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# Title
$ws.Range("B1").Value = "One2708258243"

# Labels (column A), top to bottom
$ws.Range("A2").Value = "adress"
$ws.Range("A3").Value = "country"
$ws.Range("A4").Value = "date"
$ws.Range("A5").Value = "desc"
$ws.Range("A6").Value = "name"
$ws.Range("A7").Value = "sity"

# Values (column B), top to bottom - markup stripped down to plain text
$ws.Range("B2").Value = "Adderss"
$ws.Range("B3").Value = "China"
$ws.Range("B4").Value = "23-04-2020"
$ws.Range("B5").Value = "dfgh"
$ws.Range("B6").Value = "One"
$ws.Range("B7").Value = "City"
